$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Range("Y31").Value = 793
$ws.Range("AB31").Value = 11923
$ws.Range("Y32").Value = 874
$ws.Range("AB32").Value = 13059
$ws.Range("Y33").Value = 968
$ws.Range("AB33").Value = 14363
$ws.Range("Y34").Value = 1017
$ws.Range("AB34").Value = 15195
$ws.Range("Y35").Value = 1055
$ws.Range("AB35").Value = 15792
$ws.Range("Y36").Value = 1144
$ws.Range("AB36").Value = 16905
$ws.Range("Y37").Value = 1211
$ws.Range("AB37").Value = 17839
$ws.Range("Y38").Value = 1281
$ws.Range("AB38").Value = 18862
$ws.Range("Y39").Value = 1333
$ws.Range("AB39").Value = 19904
$ws.Range("Y40").Value = 1382
$ws.Range("AB40").Value = 20832
$ws.Range("Y41").Value = 1415
$ws.Range("AB41").Value = 21397
$ws.Range("Y42").Value = 1430
$ws.Range("AB42").Value = 21780
$ws.Range("Y43").Value = 1496
$ws.Range("AB43").Value = 22522
$ws.Range("Y44").Value = 1535
$ws.Range("AB44").Value = 23172
$ws.Range("Y45").Value = 1569
$ws.Range("AB45").Value = 23842
$ws.Range("Y46").Value = 1597
$ws.Range("AB46").Value = 24510
$ws.Range("Y47").Value = 1626
$ws.Range("AB47").Value = 24959
$ws.Range("Y48").Value = 1652
$ws.Range("AB48").Value = 25422
$ws.Range("Y49").Value = 1665
$ws.Range("AB49").Value = 25694
$ws.Range("Y50").Value = 1679
$ws.Range("AB50").Value = 25945
$ws.Range("Y51").Value = 1706
$ws.Range("AB51").Value = 26267
$ws.Range("Y52").Value = 1722
$ws.Range("AB52").Value = 26585
$ws.Range("Y53").Value = 1739
$ws.Range("AB53").Value = 26884
$ws.Range("Y57").Value = 1791
$ws.Range("AB57").Value = 27871
$ws.Range("Y58").Value = 1798
$ws.Range("AB58").Value = 28030
$ws.Range("Y59").Value = 1804
$ws.Range("AB59").Value = 28226
$ws.Range("Y60").Value = 1820
$ws.Range("AB60").Value = 28425
$ws.Range("Y61").Value = 1829
$ws.Range("AB61").Value = 28590
$ws.Range("Y62").Value = 1834
$ws.Range("AB62").Value = 28713
$ws.Range("Y63").Value = 1836
$ws.Range("AB63").Value = 28775
$ws.Range("AB64").Value = 28880

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("Y58").Value = 115
$ws.Range("AB58").Value = 1536
$ws.Range("Y59").Value = 118
$ws.Range("AB59").Value = 1561
$ws.Range("Y60").Value = 123
$ws.Range("AB60").Value = 1589
$ws.Range("Y61").Value = 127
$ws.Range("AB61").Value = 1609
$ws.Range("Y62").Value = 129
$ws.Range("AB62").Value = 1626
$ws.Range("Y63").Value = 131
$ws.Range("AB63").Value = 1636
$ws.Range("Y64").Value = 132
$ws.Range("AB64").Value = 1646

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("Y33").Value = 137
$ws.Range("AB33").Value = 1793
$ws.Range("Y34").Value = 142
$ws.Range("AB34").Value = 1881
$ws.Range("Y35").Value = 152
$ws.Range("AB35").Value = 1995
$ws.Range("Y36").Value = 152
$ws.Range("AB36").Value = 2180
$ws.Range("Y37").Value = 151
$ws.Range("AB37").Value = 2208
$ws.Range("Y38").Value = 145
$ws.Range("AB38").Value = 2288
$ws.Range("Y39").Value = 145
$ws.Range("AB39").Value = 2352
$ws.Range("Y40").Value = 145
$ws.Range("AB40").Value = 2333
$ws.Range("Y41").Value = 148
$ws.Range("AB41").Value = 2321
$ws.Range("Y42").Value = 145
$ws.Range("AB42").Value = 2306
$ws.Range("Y43").Value = 132
$ws.Range("AB43").Value = 2310
$ws.Range("Y44").Value = 129
$ws.Range("AB44").Value = 2228
$ws.Range("Y45").Value = 122
$ws.Range("AB45").Value = 2140
$ws.Range("Y46").Value = 117
$ws.Range("AB46").Value = 2072
$ws.Range("Y47").Value = 115
$ws.Range("AB47").Value = 2012
$ws.Range("Y48").Value = 112
$ws.Range("AB48").Value = 1937
$ws.Range("Y49").Value = 110
$ws.Range("AB49").Value = 1914
$ws.Range("Y50").Value = 106
$ws.Range("AB50").Value = 1897
$ws.Range("Y51").Value = 100
$ws.Range("AB51").Value = 1859
$ws.Range("Y52").Value = 93
$ws.Range("AB52").Value = 1735
$ws.Range("Y53").Value = 89
$ws.Range("AB53").Value = 1679
$ws.Range("Y54").Value = 83
$ws.Range("AB54").Value = 1581
$ws.Range("Y55").Value = 82
$ws.Range("AB55").Value = 1530
$ws.Range("Y56").Value = 80
$ws.Range("AB56").Value = 1524
$ws.Range("Y57").Value = 75
$ws.Range("AB57").Value = 1500
$ws.Range("Y58").Value = 72
$ws.Range("AB58").Value = 1423
$ws.Range("Y59").Value = 73
$ws.Range("AB59").Value = 1358
$ws.Range("Y60").Value = 69
$ws.Range("AB60").Value = 1306
$ws.Range("Y61").Value = 67
$ws.Range("AB61").Value = 1260
$ws.Range("Y62").Value = 69
$ws.Range("AB62").Value = 1241
$ws.Range("Y63").Value = 69
$ws.Range("AB63").Value = 1219
$ws.Range("Y64").Value = 68
$ws.Range("AB64").Value = 1216

$ws = $wb.Worksheets.Item("ICU")
$ws.Range("Y61").Value = 14
$ws.Range("AB61").Value = 199
$ws.Range("Y62").Value = 14
$ws.Range("AB62").Value = 193
$ws.Range("Y63").Value = 15
$ws.Range("AB63").Value = 185
$ws.Range("Y64").Value = 15
$ws.Range("AB64").Value = 184

$ws = $wb.Worksheets.Item("Ventilated")
$ws.Range("Y60").Value = 10
$ws.Range("AB60").Value = 140
$ws.Range("Y64").Value = 10

$ws = $wb.Worksheets.Item("Released")
$ws.Range("Y64").Value = 220
$ws.Range("AB64").Value = 4627
